$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "scatter1"
$ws.Range("B1").Value = "scatter2"
